$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.549.76"
$ws.Range("E2").Value = "  -4.53%  "
$ws.Range("D3").Value = "2.518.74"
$ws.Range("E3").Value = "  -4.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.68"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.01"
$ws.Range("E6").Value = "  -4.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("D9").Value = "2.514.79"
$ws.Range("E9").Value = "  -4.82%  "
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").Value = "2.940.54"
$ws.Range("E14").Value = "  -5.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.50"
$ws.Range("E15").Value = "  -4.42%  "
$ws.Range("D16").Value = "59.583.29"
$ws.Range("E16").Value = "  -4.35%  "
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("D18").Value = "2.504.82"
$ws.Range("E18").Value = "  -5.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.39"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.11"
$ws.Range("E21").Value = "  -3.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.992"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  -4.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.36"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("E25").Value = "  -10.68%  "
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("E27").Value = "  -3.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.81"
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("D29").Value = "0.0₃0797"
$ws.Range("E29").Value = "  -4.43%  "
$ws.Range("E30").Value = "  -5.70%  "
$ws.Range("E31").Value = "  -4.04%  "
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.48"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.06"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("E37").Value = "  -5.20%  "
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "315.91"
$ws.Range("E40").Value = "  -5.67%  "
$ws.Range("E41").Value = "  -3.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.78"
$ws.Range("E42").Value = "  -3.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.832"
$ws.Range("E43").Value = "  -7.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.604"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "126.93"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("E50").Value = "  -1.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.72"
$ws.Range("E51").Value = "  -5.32%  "
